$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FEINmismatch")

$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Wed Nov 05 11:22:18 EST 2025"

$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Wed Nov 05 11:23:31 EST 2025"

$ws.Range("A4").Value = "Fail"
$ws.Range("B4").Value = "Wed Nov 05 11:24:43 EST 2025"

$ws.Range("A5").Value = "Fail"
$ws.Range("B5").Value = "Wed Nov 05 11:25:55 EST 2025"
